{"js": "const replacements = [\n  [\"536\u00d76=3216\", \"967\u00d79=8703\"],\n  [\"406\u00d74=1624\", \"426\u00d77=2982\"],\n  [\"555\u00d76=3330\", \"702\u00d76=4212\"],\n  [\"440\u00d79=3960\", \"190\u00d78=1520\"],\n  [\"354\u00d77=2478\", \"400\u00d73=1200\"],\n  [\"472\u00d79=4248\", \"614\u00d77=4298\"],\n  [\"251\u00d73=753\", \"722\u00d72=1444\"],\n  [\"317\u00d72=634\", \"465\u00d73=1395\"],\n  [\"317\u00d75=1585\", \"490\u00d77=3430\"],\n  [\"254\u00d79=2286\", \"403\u00d79=3627\"],\n  [\"315\u00d75=1575\", \"479\u00d74=1916\"],\n  [\"402\u00d77=2814\", \"429\u00d74=1716\"],\n  [\"957\u00d74=3828\", \"999\u00d78=7992\"],\n  [\"907\u00d79=8163\", \"388\u00d76=2328\"],\n  [\"684\u00d76=4104\", \"834\u00d76=5004\"],\n  [\"436\u00d74=1744\", \"982\u00d72=1964\"],\n  [\"805\u00d79=7245\", \"841\u00d77=5887\"],\n  [\"458\u00d74=1832\", \"916\u00d75=4580\"],\n  [\"997\u00d73=2991\", \"352\u00d73=1056\"],\n  [\"388\u00d75=1940\", \"775\u00d78=6200\"],\n  [\"851\u00d73=2553\", \"645\u00d72=1290\"],\n  [\"639\u00d77=4473\", \"933\u00d76=5598\"],\n  [\"109\u00d75=545\", \"364\u00d73=1092\"],\n  [\"528\u00d72=1056\", \"343\u00d77=2401\"],\n  [\"179\u00d73=537\", \"124\u00d75=620\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-ExactText($doc, $oldText, $newText) {\n    $r = $doc.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $null = $r.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nReplace-ExactText $d \"536\u00d76=3216\" \"967\u00d79=8703\"\nReplace-ExactText $d \"406\u00d74=1624\" \"426\u00d77=2982\"\nReplace-ExactText $d \"555\u00d76=3330\" \"702\u00d76=4212\"\nReplace-ExactText $d \"440\u00d79=3960\" \"190\u00d78=1520\"\nReplace-ExactText $d \"354\u00d77=2478\" \"400\u00d73=1200\"\nReplace-ExactText $d \"472\u00d79=4248\" \"614\u00d77=4298\"\nReplace-ExactText $d \"251\u00d73=753\" \"722\u00d72=1444\"\nReplace-ExactText $d \"317\u00d72=634\" \"465\u00d73=1395\"\nReplace-ExactText $d \"317\u00d75=1585\" \"490\u00d77=3430\"\nReplace-ExactText $d \"254\u00d79=2286\" \"403\u00d79=3627\"\nReplace-ExactText $d \"315\u00d75=1575\" \"479\u00d74=1916\"\nReplace-ExactText $d \"402\u00d77=2814\" \"429\u00d74=1716\"\nReplace-ExactText $d \"957\u00d74=3828\" \"999\u00d78=7992\"\nReplace-ExactText $d \"907\u00d79=8163\" \"388\u00d76=2328\"\nReplace-ExactText $d \"684\u00d76=4104\" \"834\u00d76=5004\"\nReplace-ExactText $d \"436\u00d74=1744\" \"982\u00d72=1964\"\nReplace-ExactText $d \"805\u00d79=7245\" \"841\u00d77=5887\"\nReplace-ExactText $d \"458\u00d74=1832\" \"916\u00d75=4580\"\nReplace-ExactText $d \"997\u00d73=2991\" \"352\u00d73=1056\"\nReplace-ExactText $d \"388\u00d75=1940\" \"775\u00d78=6200\"\nReplace-ExactText $d \"851\u00d73=2553\" \"645\u00d72=1290\"\nReplace-ExactText $d \"639\u00d77=4473\" \"933\u00d76=5598\"\nReplace-ExactText $d \"109\u00d75=545\" \"364\u00d73=1092\"\nReplace-ExactText $d \"528\u00d72=1056\" \"343\u00d77=2401\"\nReplace-ExactText $d \"179\u00d73=537\" \"124\u00d75=620\"\n"}
